# Expand preset ranges and shorten mod name
#
# Updates the "Default" (D) values for several slider rows and expands the
# shared "Values" preset lists (Duration / Cooldown / Smoothing) on both the
# "Menu Mock" sheet and the "Providers" lookup sheet.

$wb = $excel.ActiveWorkbook
$wsMenu = $wb.Worksheets.Item("Menu Mock")
$wsProviders = $wb.Worksheets.Item("Providers")

# Expanded preset value lists
$DurationValues = "0.5s | 0.6s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.5s | 1.8s | 1.875s | 2.0s | 2.25s | 2.4s | 2.5s | 3.0s | 3.75s | 4.0s | 4.5s | 5.0s | 6.0s | 6.25s | 7.5s | 10.0s"
$CooldownValues = "0s | 1.6s | 2.0s | 2.5s | 2.8s | 3.5s | 4.0s | 4.9s | 5.0s | 6.0s | 7.0s | 7.5s | 8.0s | 10.0s | 10.5s | 14.0s | 17.5s | 18.0s | 22.5s | 24.5s | 28.0s | 31.5s | 45.0s | 67.5s | 90.0s | 157.5s"
$SmoothingValues = "1.6x | 2x | 2.4x | 3x | 3.2x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 20x | 25x"

# Rows on "Menu Mock" that hold a Duration/Cooldown/Smoothing slider, with
# their (possibly updated) Default value.
$DurationRows = @(
    @{ Row = 42; Default = "0.75s" },
    @{ Row = 50; Default = "1.2s" },
    @{ Row = 58; Default = "1.2s" },
    @{ Row = 66; Default = "1.8s" },
    @{ Row = 74; Default = "2.25s" },
    @{ Row = 81; Default = "3.0s" },
    @{ Row = 88; Default = "1.0s" }
)

$CooldownRows = @(
    @{ Row = 43; Default = "4.9s" },
    @{ Row = 51; Default = "4.9s" },
    @{ Row = 59; Default = "4.9s" },
    @{ Row = 67; Default = "3.5s" },
    @{ Row = 75; Default = $null },
    @{ Row = 82; Default = "10.5s" },
    @{ Row = 89; Default = "6.0s" }
)

$SmoothingRows = @(
    @{ Row = 44; Default = "5x" },
    @{ Row = 52; Default = "5x" },
    @{ Row = 60; Default = "5x" },
    @{ Row = 68; Default = "4x" },
    @{ Row = 76; Default = "2.4x" },
    @{ Row = 83; Default = "2.4x" },
    @{ Row = 90; Default = "6x" }
)

foreach ($item in $DurationRows) {
    if ($item.Default) {
        $wsMenu.Range("D" + $item.Row).Value = $item.Default
    }
    $wsMenu.Range("E" + $item.Row).Value = $DurationValues
}

foreach ($item in $CooldownRows) {
    if ($item.Default) {
        $wsMenu.Range("D" + $item.Row).Value = $item.Default
    }
    $wsMenu.Range("E" + $item.Row).Value = $CooldownValues
}

foreach ($item in $SmoothingRows) {
    if ($item.Default) {
        $wsMenu.Range("D" + $item.Row).Value = $item.Default
    }
    $wsMenu.Range("E" + $item.Row).Value = $SmoothingValues
}

# "Providers" lookup sheet mirrors the same three preset lists.
$wsProviders.Range("B10").Value = $CooldownValues
$wsProviders.Range("B11").Value = $DurationValues
$wsProviders.Range("B12").Value = $SmoothingValues
